# Auto-generated edit script applying scheduled market-data refresh to Sheets workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 293.83334
$ws.Cells.Item(19, 10).Value = 204.33333
$ws.Cells.Item(19, 12).Value = 204.33333
$ws.Cells.Item(19, 14).Value = -554.3333299999999
$ws.Cells.Item(80, 8).Value = 734.3333
$ws.Cells.Item(80, 9).Value = 406.75
$ws.Cells.Item(80, 10).Value = 1389.5
$ws.Cells.Item(80, 11).Value = 1220.25
$ws.Cells.Item(80, 12).Value = 4168.5
$ws.Cells.Item(80, 13).Value = -222.25
$ws.Cells.Item(80, 14).Value = -6164.5
$ws.Cells.Item(83, 8).Value = 734.3333
$ws.Cells.Item(83, 9).Value = 406.75
$ws.Cells.Item(83, 10).Value = 1389.5
$ws.Cells.Item(83, 11).Value = 3660.75
$ws.Cells.Item(83, 12).Value = 12505.5
$ws.Cells.Item(83, 13).Value = 1331.25
$ws.Cells.Item(83, 14).Value = -22489.5
$ws.Cells.Item(131, 8).Value = 1161.5
$ws.Cells.Item(131, 10).Value = 2000
$ws.Cells.Item(131, 12).Value = 6000
$ws.Cells.Item(131, 14).Value = -16080
$ws.Cells.Item(132, 8).Value = 1520.5883
$ws.Cells.Item(132, 9).Value = 1155.1666
$ws.Cells.Item(132, 11).Value = 3465.4998
$ws.Cells.Item(132, 13).Value = -935.4998000000001
$ws.Cells.Item(135, 8).Value = 1395.5
$ws.Cells.Item(135, 9).Value = 1230.1111
$ws.Cells.Item(135, 11).Value = 11070.9999
$ws.Cells.Item(135, 13).Value = -8535.999900000001
$ws.Cells.Item(138, 8).Value = 1793.0454
$ws.Cells.Item(138, 9).Value = 482.33334
$ws.Cells.Item(138, 11).Value = 1447.00002
$ws.Cells.Item(138, 13).Value = 3692.99998
$ws.Cells.Item(141, 8).Value = 5295.1177
$ws.Cells.Item(141, 9).Value = 5295.1177
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 11).Value = 15885.3531
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 13).Value = -10705.3531
$ws.Cells.Item(141, 14).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 963.125
$ws.Cells.Item(4, 9).Value = 963.125
$ws.Cells.Item(4, 11).Value = 963.125
$ws.Cells.Item(4, 13).Value = -847.125
$ws.Cells.Item(61, 8).Value = 2874.75
$ws.Cells.Item(61, 9).Value = 2999.6667
$ws.Cells.Item(61, 10).Value = 2500
$ws.Cells.Item(61, 11).Value = 2999.6667
$ws.Cells.Item(61, 12).Value = 2500
$ws.Cells.Item(61, 13).Value = -2787.6667
$ws.Cells.Item(61, 14).Value = -2924
$ws.Cells.Item(74, 8).Value = 3243.4
$ws.Cells.Item(74, 9).Value = 3243.4
$ws.Cells.Item(74, 11).Value = 3243.4
$ws.Cells.Item(74, 13).Value = -2369.4
$ws.Cells.Item(77, 8).Value = 3243.4
$ws.Cells.Item(77, 9).Value = 3243.4
$ws.Cells.Item(77, 11).Value = 16217
$ws.Cells.Item(77, 13).Value = -11849
$ws.Cells.Item(132, 8).Value = 2756.7144
$ws.Cells.Item(132, 9).Value = 2756.7144
$ws.Cells.Item(132, 11).Value = 8270.143199999999
$ws.Cells.Item(132, 13).Value = -5740.143199999999
$ws.Cells.Item(136, 8).Value = 2874.75
$ws.Cells.Item(136, 9).Value = 2999.6667
$ws.Cells.Item(136, 10).Value = 2500
$ws.Cells.Item(136, 11).Value = 8999.000100000001
$ws.Cells.Item(136, 12).Value = 7500
$ws.Cells.Item(136, 13).Value = -6449.000100000001
$ws.Cells.Item(136, 14).Value = -12600

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 14966.667
$ws.Cells.Item(134, 9).Value = 14966.667
$ws.Cells.Item(134, 11).Value = 44900.001
$ws.Cells.Item(134, 13).Value = -42365.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(47, 8).Value = 50000
$ws.Cells.Item(47, 10).Value = 50000
$ws.Cells.Item(47, 12).Value = 50000
$ws.Cells.Item(47, 14).Value = -51132
$ws.Cells.Item(132, 8).Value = 1428.5
$ws.Cells.Item(132, 9).Value = 1428.5
$ws.Cells.Item(132, 11).Value = 4285.5
$ws.Cells.Item(132, 13).Value = -1755.5
$ws.Cells.Item(134, 8).Value = 5727.1665
$ws.Cells.Item(134, 9).Value = 3899
$ws.Cells.Item(134, 10).Value = 6641.25
$ws.Cells.Item(134, 11).Value = 11697
$ws.Cells.Item(134, 12).Value = 19923.75
$ws.Cells.Item(134, 13).Value = -9162
$ws.Cells.Item(134, 14).Value = -24993.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 498.75
$ws.Cells.Item(12, 9).Value = 497
$ws.Cells.Item(12, 11).Value = 1491
$ws.Cells.Item(12, 13).Value = -1318
$ws.Cells.Item(33, 8).Value = 51
$ws.Cells.Item(33, 9).Value = 16.5
$ws.Cells.Item(33, 10).Value = 120
$ws.Cells.Item(33, 11).Value = 99
$ws.Cells.Item(33, 12).Value = 720
$ws.Cells.Item(33, 13).Value = 184
$ws.Cells.Item(33, 14).Value = -1286
$ws.Cells.Item(38, 8).Value = 35
$ws.Cells.Item(38, 10).Value = 29
$ws.Cells.Item(38, 12).Value = 87
$ws.Cells.Item(38, 14).Value = -781
$ws.Cells.Item(68, 8).Value = 1499.8462
$ws.Cells.Item(68, 9).Value = 1600.25
$ws.Cells.Item(68, 10).Value = 1455.2222
$ws.Cells.Item(68, 11).Value = 4800.75
$ws.Cells.Item(68, 12).Value = 4365.6666
$ws.Cells.Item(68, 13).Value = -3989.75
$ws.Cells.Item(68, 14).Value = -5987.6666
$ws.Cells.Item(71, 8).Value = 1499.8462
$ws.Cells.Item(71, 9).Value = 1600.25
$ws.Cells.Item(71, 10).Value = 1455.2222
$ws.Cells.Item(71, 11).Value = 14402.25
$ws.Cells.Item(71, 12).Value = 13096.9998
$ws.Cells.Item(71, 13).Value = -10346.25
$ws.Cells.Item(71, 14).Value = -21208.9998
$ws.Cells.Item(86, 8).Value = 70.333336
$ws.Cells.Item(86, 10).Value = 90.666664
$ws.Cells.Item(86, 12).Value = 271.999992
$ws.Cells.Item(86, 14).Value = -2643.999992
$ws.Cells.Item(89, 8).Value = 70.333336
$ws.Cells.Item(89, 10).Value = 90.666664
$ws.Cells.Item(89, 12).Value = 815.9999759999999
$ws.Cells.Item(89, 14).Value = -12671.999976
$ws.Cells.Item(92, 8).Value = 1664.8334
$ws.Cells.Item(92, 9).Value = 1963.3334
$ws.Cells.Item(92, 10).Value = 1366.3334
$ws.Cells.Item(92, 11).Value = 5890.0002
$ws.Cells.Item(92, 12).Value = 4099.0002
$ws.Cells.Item(92, 13).Value = -4642.0002
$ws.Cells.Item(92, 14).Value = -6595.0002
$ws.Cells.Item(97, 8).Value = 1057.875
$ws.Cells.Item(97, 9).Value = 1661
$ws.Cells.Item(97, 10).Value = 454.75
$ws.Cells.Item(97, 11).Value = 4983
$ws.Cells.Item(97, 12).Value = 1364.25
$ws.Cells.Item(97, 13).Value = -4487
$ws.Cells.Item(97, 14).Value = -2356.25
$ws.Cells.Item(113, 8).Value = 1029.1666
$ws.Cells.Item(113, 9).Value = 475
$ws.Cells.Item(113, 10).Value = 1140
$ws.Cells.Item(113, 11).Value = 1425
$ws.Cells.Item(113, 12).Value = 3420
$ws.Cells.Item(113, 13).Value = 745
$ws.Cells.Item(113, 14).Value = -7760
$ws.Cells.Item(131, 8).Value = 2251.2856
$ws.Cells.Item(131, 10).Value = 2613.889
$ws.Cells.Item(131, 12).Value = 7841.667
$ws.Cells.Item(131, 14).Value = -17921.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2002.75
$ws.Cells.Item(132, 10).Value = 1014
$ws.Cells.Item(132, 12).Value = 3042
$ws.Cells.Item(132, 14).Value = -8102

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1523.1666
$ws.Cells.Item(16, 9).Value = 1567.8
$ws.Cells.Item(16, 11).Value = 1567.8
$ws.Cells.Item(16, 13).Value = -1397.8
$ws.Cells.Item(55, 8).Value = 804.3333
$ws.Cells.Item(55, 9).Value = 834.1429000000001
$ws.Cells.Item(55, 10).Value = 700
$ws.Cells.Item(55, 11).Value = 834.1429000000001
$ws.Cells.Item(55, 12).Value = 700
$ws.Cells.Item(55, 13).Value = -661.1429000000001
$ws.Cells.Item(55, 14).Value = -1046
$ws.Cells.Item(127, 8).Value = 25153.75
$ws.Cells.Item(127, 10).Value = 25153.75
$ws.Cells.Item(127, 12).Value = 25153.75
$ws.Cells.Item(127, 14).Value = -35073.75
$ws.Cells.Item(132, 8).Value = 8326.174000000001
$ws.Cells.Item(132, 9).Value = 8526.474
$ws.Cells.Item(132, 11).Value = 25579.422
$ws.Cells.Item(132, 13).Value = -23049.422

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 13).ClearContents()
$ws.Cells.Item(67, 8).Value = 0
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 11).Value = 0
$ws.Cells.Item(67, 13).ClearContents()
$ws.Cells.Item(101, 8).Value = 18614.285
$ws.Cells.Item(101, 10).Value = 18614.285
$ws.Cells.Item(101, 12).Value = 18614.285
$ws.Cells.Item(101, 14).Value = -25104.285
$ws.Cells.Item(126, 8).Value = 4322.3335
$ws.Cells.Item(126, 9).Value = 4454.625
$ws.Cells.Item(126, 10).Value = 4171.143
$ws.Cells.Item(126, 11).Value = 13363.875
$ws.Cells.Item(126, 12).Value = 12513.429
$ws.Cells.Item(126, 13).Value = -10893.875
$ws.Cells.Item(126, 14).Value = -17453.429

Write-Host "Applied 194 cell updates across 8 sheets"